$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# The "Ready for handoff" status text is shared (as the same backing string)
# by Overview!B2, Overview!C2, zh-cn!B2 and de-de!B2. Update all of them so
# they keep pointing at one shared piece of text, now reading
# "Handoff transform failed".
$wsOverview.Range("B2").Value = "Handoff transform failed"
$wsOverview.Range("C2").Value = "Handoff transform failed"

# --- zh-cn sheet ---
# The collection-level Hyperlinks.Delete() removes every hyperlink on the
# sheet, so capture what needs to survive first, wipe them all, then restore
# the ones that are not being dropped (A2 / A3). C2's hyperlink (the stale
# handoff-file link) is intentionally not recreated.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0c17116158ede433807cbb0d129d0a6715d13f0b/e2e/351eeb8b-9243-4783-8014-21db8f5b9bc8.md", "", "", "351eeb8b-9243-4783-8014-21db8f5b9bc8.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0c17116158ede433807cbb0d129d0a6715d13f0b/.localization-config", "", "", ".localization-config")

$wsZh.Range("B2").Value = "Handoff transform failed"
$wsZh.Range("C2").Clear()
$wsZh.Range("D2").Value = "0001-01-01 00:00:00"
$wsZh.Range("G2").Value = "0001-01-01 00:00:00"
$wsZh.Range("H2").Value = "Ignored"

# --- de-de sheet ---
$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/0c17116158ede433807cbb0d129d0a6715d13f0b/e2e/351eeb8b-9243-4783-8014-21db8f5b9bc8.md", "", "", "351eeb8b-9243-4783-8014-21db8f5b9bc8.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/0c17116158ede433807cbb0d129d0a6715d13f0b/.localization-config", "", "", ".localization-config")

$wsDe.Range("B2").Value = "Handoff transform failed"
$wsDe.Range("C2").Clear()
$wsDe.Range("D2").Value = "0001-01-01 00:00:00"
$wsDe.Range("G2").Value = "0001-01-01 00:00:00"
$wsDe.Range("H2").Value = "Ignored"
